# Updated loading_percent results for the 380 kV case (Case_3_72)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> column letter -> new value
$newValues = @{
    2 = @{ "B" = 16.02012629397998; "C" = 9.681425524113124; "D" = 5.484425448639597; "F" = 33.000050942148; "G" = 3.638034897653364; "M" = 19.97039798536258; "N" = 18.44959073503417 }
    3 = @{ "B" = 15.43027468497448; "C" = 9.114404059122982; "D" = 5.50022711893808; "F" = 32.48786273960563; "G" = 3.642272097676737; "M" = 19.36962997078813; "N" = 18.49451648002974 }
    4 = @{ "B" = 15.06184701594167; "C" = 8.75012428297048; "D" = 5.510772248401409; "F" = 32.18032167991367; "G" = 3.645004664689208; "M" = 18.99957655292726; "N" = 18.52403160506164 }
    5 = @{ "B" = 14.91040443096195; "C" = 8.597766091800199; "D" = 5.515280204693637; "F" = 32.05689670160166; "G" = 3.646151268598247; "M" = 18.84873168030756; "N" = 18.53654326304103 }
    6 = @{ "B" = 14.88518695861828; "C" = 8.572235707831965; "D" = 5.516041440048851; "F" = 32.036521285918; "G" = 3.646343662262307; "M" = 18.82368886195226; "N" = 18.53865000537662 }
    7 = @{ "B" = 15.05980952205379; "C" = 8.74808515118732; "D" = 5.510832192610712; "F" = 32.17864923174373; "G" = 3.645019994138737; "M" = 18.99754204767037; "N" = 18.52419838380234 }
    8 = @{ "B" = 15.81819734991304; "C" = 9.489340897375982; "D" = 5.489698200081642; "F" = 32.82210082673336; "G" = 3.639468800099103; "M" = 19.76364493340837; "N" = 18.46467957281216 }
    9 = @{ "B" = 17.24519560674032; "C" = 10.81004017571188; "D" = 5.454997178196623; "F" = 34.13207192719405; "G" = 3.62961511530608; "M" = 21.24625000321489; "N" = 18.36334681983316 }
    10 = @{ "B" = 18.24445958242246; "C" = 11.75688686651017; "D" = 5.433690039821119; "F" = 35.11447971230285; "G" = 3.622995723795056; "M" = 22.31072201699047; "N" = 18.2983671845226 }
    11 = @{ "B" = 18.68620657229689; "C" = 12.17092546747112; "D" = 5.42492335045868; "F" = 35.56375834955728; "G" = 3.620117078203849; "M" = 22.78722912428322; "N" = 18.27088307220303 }
    12 = @{ "B" = 18.85148846640055; "C" = 12.32414781487357; "D" = 5.421738264450293; "F" = 35.73406840316685; "G" = 3.619045917796649; "M" = 22.96638549605407; "N" = 18.26077583522405 }
    13 = @{ "B" = 18.81598316722958; "C" = 12.2913068998424; "D" = 5.422418215935421; "F" = 35.69738372892542; "G" = 3.619275772107321; "M" = 22.92786079326175; "N" = 18.26293921744784 }
    14 = @{ "B" = 18.6998452718349; "C" = 12.18360250423844; "D" = 5.424658603812379; "F" = 35.57776721503161; "G" = 3.620028574807051; "M" = 22.80199515381687; "N" = 18.27004551190236 }
    15 = @{ "B" = 18.62844300426433; "C" = 12.11716673003486; "D" = 5.426048489188655; "F" = 35.50451704147103; "G" = 3.620492148004077; "M" = 22.7247263617147; "N" = 18.27443750224603 }
    16 = @{ "B" = 18.21531871481193; "C" = 11.72932856650982; "D" = 5.434281725819941; "F" = 35.08515202114884; "G" = 3.623186505296176; "M" = 22.2794102810647; "N" = 18.30020526490939 }
    17 = @{ "B" = 17.95848002141884; "C" = 11.48503248838876; "D" = 5.439570791324742; "F" = 34.82837599621398; "G" = 3.624873255450454; "M" = 22.00411393581392; "N" = 18.31654589697537 }
    18 = @{ "B" = 17.80955412112034; "C" = 11.34217290967379; "D" = 5.442699992702877; "F" = 34.68091846124826; "G" = 3.62585591403404; "M" = 21.84505287039906; "N" = 18.32613982743214 }
    19 = @{ "B" = 17.75892944910617; "C" = 11.29339945439154; "D" = 5.443774395602418; "F" = 34.63103699156991; "G" = 3.626190773839286; "M" = 21.79107986160227; "N" = 18.32942163987922 }
    20 = @{ "B" = 17.9859462177323; "C" = 11.51128098662849; "D" = 5.438998738969223; "F" = 34.85568730257533; "G" = 3.624692406932978; "M" = 22.0334953269998; "N" = 18.31478618916305 }
    21 = @{ "B" = 18.73401314675089; "C" = 12.21533450172105; "D" = 5.423996880542985; "F" = 35.61289793691773; "G" = 3.619806946034167; "M" = 22.83900117110933; "N" = 18.26795005234732 }
    22 = @{ "B" = 19.21121356847899; "C" = 12.65470288112219; "D" = 5.414977974096889; "F" = 36.10874372351228; "G" = 3.616724231224686; "M" = 23.35787774947529; "N" = 18.23909183535036 }
    23 = @{ "B" = 18.95763965231743; "C" = 12.42209807406518; "D" = 5.419719138967178; "F" = 35.84406609621683; "G" = 3.618359495666063; "M" = 23.08168951680996; "N" = 18.25433304964977 }
    24 = @{ "B" = 17.97353268151613; "C" = 11.49942155018848; "D" = 5.439257088459383; "F" = 34.84333933013057; "G" = 3.624774128283849; "M" = 22.02021444574808; "N" = 18.31558113186606 }
    25 = @{ "B" = 16.86697603195838; "C" = 10.46758108519779; "D" = 5.463655010124763; "F" = 33.77350310296548; "G" = 3.632171239987573; "M" = 20.84866027948991; "N" = 18.3891038558883 }
}

foreach ($row in $newValues.Keys) {
    $rowData = $newValues[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
